# Insert a new data row at row 412 (pushing the existing rows 412:478
# down to 413:479) and populate it with a new Coliflor price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("412").Insert()

$ws.Cells.Item(412, 1).Value = 7
$ws.Cells.Item(412, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(412, 3).Value = "Ñuble"
$ws.Cells.Item(412, 4).Value = 45077
$ws.Cells.Item(412, 5).Value = 16
$ws.Cells.Item(412, 6).Value = 100112008
$ws.Cells.Item(412, 7).Value = "Coliflor"
$ws.Cells.Item(412, 8).Value = "Sin especificar"
$ws.Cells.Item(412, 9).Value = "Primera"
$ws.Cells.Item(412, 10).Value = 300
$ws.Cells.Item(412, 11).Value = 1200
$ws.Cells.Item(412, 12).Value = 1200
$ws.Cells.Item(412, 13).Value = 1200
$ws.Cells.Item(412, 14).Value = "$/unidad"
$ws.Cells.Item(412, 15).Value = "Región del Maule"
$ws.Cells.Item(412, 16).Value = 1200
$ws.Cells.Item(412, 17).Value = 1
$ws.Cells.Item(412, 18).Value = "Hortaliza"
